# "Update countries & provincias Spain"
#
# The underlying COVID table is keyed by country name in column A with
# stats in columns B:H. The source data feed re-ordered a few countries
# in its source list (Peru/Paises Bajos, Benin, Namibia) and refreshed
# the case counts for this snapshot. Because the sheet has no blank
# rows to insert/delete around, the net effect of the re-ordering is
# just that a handful of rows swap/shift which country (and therefore
# which stats) they hold, while every other row is untouched. We apply
# the resulting cell-level deltas directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 20:52"

# Estados Unidos (row 4) - refreshed totals
$ws.Range("B4").Value = 1112771
$ws.Range("C4").Value = 17748
$ws.Range("D4").Value = 158029
$ws.Range("E4").Value = 889811
$ws.Range("G4").Value = 1075
$ws.Range("H4").Value = 64931

# Francia (row 8) - refreshed totals
$ws.Range("B8").Value = 167346
$ws.Range("C8").Value = 168
$ws.Range("E8").Value = 92540

# Canada (row 15) - refreshed totals
$ws.Range("B15").Value = 54784
$ws.Range("C15").Value = 1548
$ws.Range("E15").Value = 29302
$ws.Range("G15").Value = 203
$ws.Range("H15").Value = 3387

# Peru moves above Paises Bajos (rows 17/18 swap country + stats);
# Peru also gets refreshed totals, Paises Bajos keeps its prior totals
$ws.Range("A17").Value = "Peru"
$ws.Range("B17").Value = 40459
$ws.Range("C17").Value = 3483
$ws.Range("D17").Value = 11129
$ws.Range("E17").Value = 28206
$ws.Range("F17").Value = 658
$ws.Range("G17").Value = 73
$ws.Range("H17").Value = 1124

$ws.Range("A18").Value = "Paises Bajos"
$ws.Range("B18").Value = 39791
$ws.Range("C18").Value = 475
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 34648
$ws.Range("F18").Value = 735
$ws.Range("G18").Value = 98
$ws.Range("H18").Value = 4893

# Malta (row 115) - refreshed totals
$ws.Range("D115").Value = 383
$ws.Range("E115").Value = 80

# Benin moves up to right after Monaco, pushing Uganda..Republica de
# Africa Central (rows 156-166) down by one row each; Benin itself
# gets refreshed totals while every displaced country keeps carrying
# its own previous totals down to its new row.
$ws.Range("A155").Value = "Benin"
$ws.Range("B155").Value = 90
$ws.Range("C155").Value = 26
$ws.Range("D155").Value = 42
$ws.Range("E155").Value = 46
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 2

$ws.Range("A156").Value = "Uganda"
$ws.Range("B156").Value = 83
$ws.Range("D156").Value = 52
$ws.Range("E156").Value = 31
$ws.Range("F156").Value = 0
$ws.Range("H156").Value = 0

$ws.Range("A157").Value = "Guyana"
$ws.Range("D157").Value = 22
$ws.Range("E157").Value = 51
$ws.Range("F157").Value = 2
$ws.Range("H157").Value = 9

$ws.Range("A158").Value = "Liechtenstein"
$ws.Range("B158").Value = 82
$ws.Range("D158").Value = 55
$ws.Range("E158").Value = 26
$ws.Range("H158").Value = 1

$ws.Range("A159").Value = "Haiti"
$ws.Range("D159").Value = 8
$ws.Range("E159").Value = 65
$ws.Range("F159").Value = 0
$ws.Range("H159").Value = 8

$ws.Range("A160").Value = "Bahamas"
$ws.Range("D160").Value = 25
$ws.Range("E160").Value = 45
$ws.Range("F160").Value = 1
$ws.Range("H160").Value = 11

$ws.Range("A161").Value = "Barbados"
$ws.Range("B161").Value = 81
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 39
$ws.Range("E161").Value = 35
$ws.Range("F161").Value = 4
$ws.Range("H161").Value = 7

$ws.Range("A162").Value = "Mozambique"
$ws.Range("B162").Value = 79
$ws.Range("C162").Value = 3
$ws.Range("D162").Value = 12
$ws.Range("E162").Value = 67
$ws.Range("F162").Value = 0
$ws.Range("H162").Value = 0

$ws.Range("A163").Value = "San Martin (Parte Holandesa)"
$ws.Range("B163").Value = 76
$ws.Range("C163").Value = 1
$ws.Range("D163").Value = 44
$ws.Range("E163").Value = 19
$ws.Range("F163").Value = 7
$ws.Range("H163").Value = 13

$ws.Range("A164").Value = "Islas Caimanes"
$ws.Range("D164").Value = 10
$ws.Range("E164").Value = 62
$ws.Range("F164").Value = 3
$ws.Range("H164").Value = 1

$ws.Range("A165").Value = "Republica del Chad"
$ws.Range("B165").Value = 73
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 33
$ws.Range("E165").Value = 35
$ws.Range("H165").Value = 5

$ws.Range("A166").Value = "Republica de Africa Central"
$ws.Range("C166").Value = 14
$ws.Range("D166").Value = 10
$ws.Range("E166").Value = 54
$ws.Range("H166").Value = 0

# Angola (row 180) - refreshed totals
$ws.Range("B180").Value = 30
$ws.Range("C180").Value = 3
$ws.Range("D180").Value = 11
$ws.Range("E180").Value = 17

# Botsuana (row 183) - refreshed totals
$ws.Range("D183").Value = 8
$ws.Range("E183").Value = 14

# Namibia moves above San Vicente y las Granadinas (rows 192/193 swap
# country names only - both rows already shared identical stats)
$ws.Range("A192").Value = "Namibia"
$ws.Range("A193").Value = "San Vicente y las Granadinas"
